$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of results (row 9) with the new hyperparameters string,
# train accuracy and test accuracy - continuing the existing results table.
$ws.Range("A9").Value = "{'criterion': 'entropy', 'max_depth': 7, 'max_features': 'sqrt', 'max_leaf_nodes': 10, 'min_samples_leaf': 3, 'min_samples_split': 6, 'n_estimators': 2000}"
$ws.Range("B9").Value = 0.896
$ws.Range("C9").Value = 0.8

# Match the formatting of the other data rows (no explicit per-cell style)
# rather than the implicit column-default style new cells pick up.
$ws.Range("A2:C2").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)

